$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated crypto price / 1h-volume data, plus the Quant <-> PaxDollar row swap ---

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "310.51", "20.40", "0.000008746") must be forced to Text format first
# so the literal digits/trailing zeros in the source data are preserved.
$textForcedUpdates = @{
    "D5" = "310.51"
    "D9" = "0.07383"
    "D10" = "0.8741"
    "D11" = "20.40"
    "D13" = "5.354"
    "D14" = "6.542"
    "D15" = "0.07056"
    "D16" = "91.25"
    "D18" = "0.000008746"
    "D20" = "14.77"
    "D22" = "5.314"
    "D23" = "10.79"
    "D25" = "1.916"
    "D26" = "151.39"
    "D27" = "18.51"
    "D28" = "2.153"
    "D29" = "5.295"
    "D30" = "116.04"
    "D31" = "0.08901"
    "D32" = "0.7719"
    "D35" = "2.908"
    "D36" = "1.001"
    "D38" = "0.01958"
    "D39" = "0.05242"
    "D40" = "2.416"
    "D41" = "0.5382"
    "D42" = "7.263"
    "D43" = "2.905"
    "D44" = "0.1659"
    "D45" = "8.556"
    "D46" = "0.5060"
    "D47" = "10.35"
    "D48" = "1.001"
    "D49" = "103.99"
    "D50" = "1.662"
    "D51" = "0.06322"
}

foreach ($cellRef in $textForcedUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $textForcedUpdates[$cellRef]
}

# Remaining cells (names, links, percentage-change text, and multi-dot prices
# like "26.883.91") are unambiguous text as-is, so a plain assignment keeps them
# as normal (unstyled) strings exactly like the rest of the sheet.
$plainUpdates = @{
    "D2" = "26.883.91"
    "E2" = "  -1.90%  "
    "D3" = "1.808.44"
    "E3" = "  -0.79%  "
    "E4" = "  +0.11%  "
    "E5" = "  -1.15%  "
    "E6" = "  +0.11%  "
    "E7" = "  +2.79%  "
    "E8" = "  -0.39%  "
    "E9" = "  -1.60%  "
    "E11" = "  -2.97%  "
    "D12" = "1.770.84"
    "E12" = "  -2.83%  "
    "E13" = "  -1.20%  "
    "E14" = "  -3.27%  "
    "E15" = "  -0.63%  "
    "E16" = "  -2.90%  "
    "E18" = "  -0.64%  "
    "E19" = "  +0.08%  "
    "E20" = "  -2.86%  "
    "D21" = "26.898.12"
    "E21" = "  -1.82%  "
    "E22" = "  +0.67%  "
    "E23" = "  -1.30%  "
    "D24" = "2.004.63"
    "E24" = "  -2.50%  "
    "E25" = "  -2.04%  "
    "E26" = "  -0.10%  "
    "E27" = "  -0.37%  "
    "E28" = "  -9.21%  "
    "E29" = "  -1.55%  "
    "E30" = "  -1.72%  "
    "E31" = "  +0.58%  "
    "E32" = "  -2.02%  "
    "E33" = "  -3.32%  "
    "E34" = "  -0.87%  "
    "E35" = "  -0.47%  "
    "E36" = "  +0.11%  "
    "E37" = "  +0.27%  "
    "E38" = "  -1.96%  "
    "E39" = "  -1.58%  "
    "E40" = "  +4.59%  "
    "E41" = "  +0.99%  "
    "E42" = "  -1.38%  "
    "E43" = "  +1.77%  "
    "E44" = "  -3.86%  "
    "E45" = "  -2.17%  "
    "E46" = "  -0.70%  "
    "E47" = "  -2.56%  "
    "B48" = "PaxDollar"
    "C48" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "E48" = "  +0.10%  "
    "B49" = "Quant"
    "C49" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "E49" = "  -1.72%  "
    "E50" = "  -2.43%  "
    "E51" = "  -0.94%  "
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}
